# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - reuse the header style/format already used by G1 ("sum")
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data values for H2:H9 taken from the diff
$saveValues = @(0, 0, 0, 1, 0, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
